$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("largepart")
$ws.Range("C2").Value = "./datafiles/pig/안심살1.png"
$ws.Range("C3").Value = "./datafiles/pig/등심살1.png"
$ws.Range("C4").Value = "./datafiles/pig/알등심살1.png"
$ws.Range("C5").Value = "./datafiles/pig/등심덧살1.png"
$ws.Range("C6").Value = "./datafiles/pig/목심살1.png"
$ws.Range("C7").Value = "./datafiles/pig/앞다리살1.png"
$ws.Range("C8").Value = "./datafiles/pig/앞사태살1.png"
$ws.Range("C9").Value = "./datafiles/pig/항정살1.png"
$ws.Range("C10").Value = "./datafiles/pig/꾸리살1.png"
$ws.Range("C11").Value = "./datafiles/pig/부채살1.png"
$ws.Range("C12").Value = "./datafiles/pig/주걱살1.png"
$ws.Range("C13").Value = "./datafiles/pig/볼기살1.png"
$ws.Range("C14").Value = "./datafiles/pig/설깃살1.png"
$ws.Range("C15").Value = "./datafiles/pig/도가니살1.png"
$ws.Range("C16").Value = "./datafiles/pig/홍두깨살1.png"
$ws.Range("C17").Value = "./datafiles/pig/보섭살1.png"
$ws.Range("C18").Value = "./datafiles/pig/뒷사태살1.png"
$ws.Range("C19").Value = "./datafiles/pig/삼겹살1.png"
$ws.Range("C20").Value = "./datafiles/pig/갈매기살1.png"
$ws.Range("C21").Value = "./datafiles/pig/등갈비1.png"
$ws.Range("C22").Value = "./datafiles/pig/토시살1.png"
$ws.Range("C23").Value = "./datafiles/pig/오돌삼겹1.png"
$ws.Range("C24").Value = "./datafiles/pig/갈비1.png"
$ws.Range("C25").Value = "./datafiles/pig/갈비살1.png"
$ws.Range("C26").Value = "./datafiles/pig/마구리1.png"
$ws.Range("C27").Value = "./datafiles/cow/안심살1.png"
$ws.Range("C28").Value = "./datafiles/cow/윗등심살1.png"
$ws.Range("C29").Value = "./datafiles/cow/꽃등심살1.png"
$ws.Range("C30").Value = "./datafiles/cow/아래등심살1.png"
$ws.Range("C31").Value = "./datafiles/cow/살치살1.png"
$ws.Range("C32").Value = "./datafiles/cow/채끝살1.png"
$ws.Range("C33").Value = "./datafiles/cow/목심살1.png"
$ws.Range("C34").Value = "./datafiles/cow/꾸리살1.png"
$ws.Range("C35").Value = "./datafiles/cow/부채살1.png"
$ws.Range("C36").Value = "./datafiles/cow/앞다리살1.png"
$ws.Range("C37").Value = "./datafiles/cow/갈비덧살1.png"
$ws.Range("C38").Value = "./datafiles/cow/부채덮개살1.png"
$ws.Range("C39").Value = "./datafiles/cow/우둔살1.png"
$ws.Range("C40").Value = "./datafiles/cow/홍두깨살1.png"
$ws.Range("C41").Value = "./datafiles/cow/보섭살1.png"
$ws.Range("C42").Value = "./datafiles/cow/설깃살1.png"
$ws.Range("C43").Value = "./datafiles/cow/설깃머리살1.png"
$ws.Range("C44").Value = "./datafiles/cow/도가니살1.png"
$ws.Range("C45").Value = "./datafiles/cow/삼각살1.png"
$ws.Range("C46").Value = "./datafiles/cow/양지머리1.png"
$ws.Range("C47").Value = "./datafiles/cow/차돌박이1.png"
$ws.Range("C48").Value = "./datafiles/cow/업진살1.png"
$ws.Range("C49").Value = "./datafiles/cow/업진안살1.png"
$ws.Range("C50").Value = "./datafiles/cow/치마양지1.png"
$ws.Range("C51").Value = "./datafiles/cow/치마살1.png"
$ws.Range("C52").Value = "./datafiles/cow/앞치마살1.png"
$ws.Range("C53").Value = "./datafiles/cow/앞사태1.png"
$ws.Range("C54").Value = "./datafiles/cow/뒷사태1.png"
$ws.Range("C55").Value = "./datafiles/cow/뭉치사태1.png"
$ws.Range("C56").Value = "./datafiles/cow/아롱사태1.png"
$ws.Range("C57").Value = "./datafiles/cow/상박살1.png"
$ws.Range("C58").Value = "./datafiles/cow/본갈비1.png"
$ws.Range("C59").Value = "./datafiles/cow/꽃갈비1.png"
$ws.Range("C60").Value = "./datafiles/cow/참갈비1.png"
$ws.Range("C61").Value = "./datafiles/cow/갈비살1.png"
$ws.Range("C62").Value = "./datafiles/cow/마구리1.png"
$ws.Range("C63").Value = "./datafiles/cow/토시살1.png"
$ws.Range("C64").Value = "./datafiles/cow/안창살1.png"
$ws.Range("C65").Value = "./datafiles/cow/제비추리1.png"
$ws.Range("C23").Select()
